$wb = $excel.ActiveWorkbook

$xlShiftToRight = -4161
$xlShiftToLeft  = -4159

# --- Sheets "檢核表" and "應收應付" ---
# Insert 4 new header columns (前日餘額, 加, 減, 淨增減) right after column C
# (會計帳餘額), pushing the existing 銷帳檔餘額 / 主檔餘額 / 會計檔與主檔差額 /
# 銷帳檔與主檔差額 headers from D:G over to H:K.
foreach ($name in @("檢核表", "應收應付")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D1:G1").Insert($xlShiftToRight)
    $ws.Range("D1").Value = "前日餘額"
    $ws.Range("E1").Value = "加"
    $ws.Range("F1").Value = "減"
    $ws.Range("G1").Value = "淨增減"
}

# --- Sheet "未銷帳" ---
# Remove the 會計帳餘額 column (old C) and the 會計檔與主檔差額 column (old F),
# shifting the remaining header cells left so only 區隔帳冊 / 科目 / 銷帳檔餘額 /
# 主檔餘額 / 銷帳檔與主檔差額 remain (A:E).
$ws3 = $wb.Worksheets.Item("未銷帳")
$ws3.Range("C1").Delete($xlShiftToLeft)
$ws3.Range("E1").Delete($xlShiftToLeft)

# --- Restore the per-sheet selections left behind by the editor, ending on
# the "檢核表" tab so it stays the active sheet (matches tabSelected="1"). ---
$wb.Worksheets.Item("應收應付").Range("E13").Select() | Out-Null
$wb.Worksheets.Item("未銷帳").Range("D7").Select() | Out-Null
$wb.Worksheets.Item("檢核表").Range("E18").Select() | Out-Null
